$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 740
$ws1.Range("F5").Value = 2972
$ws1.Range("F6").Value = 60
$ws1.Range("F7").Value = 3820
$ws1.Range("F8").Value = 473
$ws1.Range("F9").Value = 956
$ws1.Range("F10").Value = 24

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 740
$ws4.Range("F6").Value = 2972
$ws4.Range("F7").Value = 60
$ws4.Range("F8").Value = 3820
$ws4.Range("F9").Value = 473
$ws4.Range("F10").Value = 956
$ws4.Range("F11").Value = 24
